$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (B1) on both sheets
$wsInput.Range("B1").Value = "2530-MS-EPP-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-1-DATE-VAR-INST-1st"
$wsOutput.Range("B1").Value = "2530-MS-EPP-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-1-DATE-VAR-INST-1st"

# Update short name (B2) on input sheet - now text instead of a number
$wsInput.Range("B2").Value = "253d"

# Move selection on input sheet to B1
$wsInput.Range("B1").Select()
